$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the updated cells keep their original Text storage type (matching the
# source workbook, where Price/Volume/Hora columns are stored as text strings,
# not numbers) by forcing a Text number format before writing the new values.
$updatedCells = @("D2","E2","G2","D3","E3","G3","D4","E4","G4","D5","E5","G5","D6","E6","G6","E7","G7","D8","E8","G8","E9","G9","D10","E10","G10","D11","E11","G11","D12","E12","G12","D13","E13","G13","E14","G14","D15","E15","G15","D16","E16","G16","G17","D18","E18","G18","E19","G19","E20","G20","D21","E21","G21","D22","E22","G22","E23","G23","D24","E24","G24","E25","G25","D26","E26","G26","D27","E27","G27","G28","G29","G30","G31","G32","G33","G34","G35","G36","G37","G38","D39","E39","G39","D40","E40","G40","D41","E41","G41","D42","G42","E43","G43","D44","E44","G44","D45","E45","G45","D46","E46","G46","E47","G47","G48","D49","E49","G49","D50","E50","G50","D51","E51","G51")
foreach ($addr in $updatedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "291.45"
$ws.Range("E2").Value = "-2.96%"
$ws.Range("G2").Value = "9"
$ws.Range("D3").Value = "30.67"
$ws.Range("E3").Value = "-6.32%"
$ws.Range("G3").Value = "9"
$ws.Range("D4").Value = "4.953"
$ws.Range("E4").Value = "0.41%"
$ws.Range("G4").Value = "9"
$ws.Range("D5").Value = "0.07216"
$ws.Range("E5").Value = "-6.43%"
$ws.Range("G5").Value = "9"
$ws.Range("D6").Value = "1.851"
$ws.Range("E6").Value = "-6.48%"
$ws.Range("G6").Value = "9"
$ws.Range("E7").Value = "-1.89%"
$ws.Range("G7").Value = "9"
$ws.Range("D8").Value = "3.768"
$ws.Range("E8").Value = "-0.85%"
$ws.Range("G8").Value = "9"
$ws.Range("E9").Value = "-2.49%"
$ws.Range("G9").Value = "9"
$ws.Range("D10").Value = "0.1657"
$ws.Range("E10").Value = "-5.94%"
$ws.Range("G10").Value = "9"
$ws.Range("D11").Value = "0.07715"
$ws.Range("E11").Value = "-0.66%"
$ws.Range("G11").Value = "9"
$ws.Range("D12").Value = "0.07995"
$ws.Range("E12").Value = "-6.89%"
$ws.Range("G12").Value = "9"
$ws.Range("D13").Value = "0.03040"
$ws.Range("E13").Value = "-3.92%"
$ws.Range("G13").Value = "9"
$ws.Range("E14").Value = "-0.03%"
$ws.Range("G14").Value = "9"
$ws.Range("D15").Value = "0.001493"
$ws.Range("E15").Value = "-1.25%"
$ws.Range("G15").Value = "9"
$ws.Range("D16").Value = "0.005740"
$ws.Range("E16").Value = "-1.67%"
$ws.Range("G16").Value = "9"
$ws.Range("G17").Value = "9"
$ws.Range("D18").Value = "3.472"
$ws.Range("E18").Value = "0.35%"
$ws.Range("G18").Value = "9"
$ws.Range("E19").Value = "-3.26%"
$ws.Range("G19").Value = "9"
$ws.Range("E20").Value = "-0.83%"
$ws.Range("G20").Value = "9"
$ws.Range("D21").Value = "0.1299"
$ws.Range("E21").Value = "-2.04%"
$ws.Range("G21").Value = "9"
$ws.Range("D22").Value = "4.035"
$ws.Range("E22").Value = "-6.76%"
$ws.Range("G22").Value = "9"
$ws.Range("E23").Value = "13.05%"
$ws.Range("G23").Value = "9"
$ws.Range("D24").Value = "0.04498"
$ws.Range("E24").Value = "-1.10%"
$ws.Range("G24").Value = "9"
$ws.Range("E25").Value = "-0.76%"
$ws.Range("G25").Value = "9"
$ws.Range("D26").Value = "0.004010"
$ws.Range("E26").Value = "-9.08%"
$ws.Range("G26").Value = "9"
$ws.Range("D27").Value = "0.0001201"
$ws.Range("E27").Value = "-4.06%"
$ws.Range("G27").Value = "9"
$ws.Range("G28").Value = "9"
$ws.Range("G29").Value = "9"
$ws.Range("G30").Value = "9"
$ws.Range("G31").Value = "9"
$ws.Range("G32").Value = "9"
$ws.Range("G33").Value = "9"
$ws.Range("G34").Value = "9"
$ws.Range("G35").Value = "9"
$ws.Range("G36").Value = "9"
$ws.Range("G37").Value = "9"
$ws.Range("G38").Value = "9"
$ws.Range("D39").Value = "0.01600"
$ws.Range("E39").Value = "-5.77%"
$ws.Range("G39").Value = "9"
$ws.Range("D40").Value = "0.04412"
$ws.Range("E40").Value = "-5.42%"
$ws.Range("G40").Value = "9"
$ws.Range("D41").Value = "0.007274"
$ws.Range("E41").Value = "-4.06%"
$ws.Range("G41").Value = "9"
$ws.Range("D42").Value = "0.009921"
$ws.Range("G42").Value = "9"
$ws.Range("E43").Value = "-3.16%"
$ws.Range("G43").Value = "9"
$ws.Range("D44").Value = "0.002052"
$ws.Range("E44").Value = "-12.16%"
$ws.Range("G44").Value = "9"
$ws.Range("D45").Value = "0.009518"
$ws.Range("E45").Value = "-16.82%"
$ws.Range("G45").Value = "9"
$ws.Range("D46").Value = "0.00005940"
$ws.Range("E46").Value = "-5.24%"
$ws.Range("G46").Value = "9"
$ws.Range("E47").Value = "-0.11%"
$ws.Range("G47").Value = "9"
$ws.Range("G48").Value = "9"
$ws.Range("D49").Value = "0.003004"
$ws.Range("E49").Value = "-3.26%"
$ws.Range("G49").Value = "9"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").Value = "-0.11%"
$ws.Range("G50").Value = "9"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").Value = "-0.11%"
$ws.Range("G51").Value = "9"
